$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 and row 5 each get a new "Expected Title" column inserted right after
# column A. Only these two rows are affected (rows 1,2,7,8 keep their existing
# layout), so we shift the existing B..G contents of rows 4 and 5 one column to
# the right (copying right-to-left so we don't clobber data before reading it),
# and then fill the now-empty column B (and the updated C5) with the new data.
foreach ($r in 4, 5) {
    for ($c = 7; $c -ge 2; $c--) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r, $c + 1)
        $src.Copy($dst)
    }
}

$ws.Range("C5").Value2 = "Ravi0111"
$ws.Range("B4").Value2 = "Expected Title"
$ws.Range("B5").Value2 = "Adactin.com - New User Registration"

# The hyperlink that used to live on F5 now belongs on G5 (it followed the
# Email Address column during the shift above).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:ravi12@gmail.com") | Out-Null

# Column B/C widths were widened to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 38.17
$ws.Columns.Item(3).ColumnWidth = 14.33

# The active selection ends up on the newly inserted H4 cell.
$ws.Range("H4").Select() | Out-Null
